# Auto-generated edit script: update cryptos list values (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.157.89"
$ws.Range("E2").Value = "  +3.52%  "
$ws.Range("D3").Value = "2.476.36"
$ws.Range("E3").Value = "  +2.59%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'488.28"
$ws.Range("E5").Value = "  +5.07%  "
$ws.Range("D6").Value = "'146.82"
$ws.Range("E6").Value = "  +11.66%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +3.74%  "
$ws.Range("D9").Value = "2.487.61"
$ws.Range("E9").Value = "  +2.38%  "
$ws.Range("D10").Value = "'5.81"
$ws.Range("E10").Value = "  +9.44%  "
$ws.Range("D11").Value = "'0.0969"
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("D12").Value = "'0.332"
$ws.Range("E12").Value = "  +5.82%  "
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").Value = "2.915.17"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("D15").Value = "56.186.26"
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").Value = "'21.11"
$ws.Range("E16").Value = "  +7.29%  "
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("D18").Value = "2.488.11"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").Value = "'4.52"
$ws.Range("E19").Value = "  +8.27%  "
$ws.Range("D20").Value = "'10.06"
$ws.Range("E20").Value = "  +6.68%  "
$ws.Range("D21").Value = "'318.05"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'5.80"
$ws.Range("E23").Value = "  +8.23%  "
$ws.Range("D24").Value = "'58.44"
$ws.Range("E24").Value = "  +4.18%  "
$ws.Range("D25").Value = "'0.412"
$ws.Range("E25").Value = "  +7.27%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("E27").Value = "  +4.29%  "
$ws.Range("D28").Value = "2.584.16"
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("E29").Value = "  +7.22%  "
$ws.Range("D30").Value = "0.0₃0789"
$ws.Range("E30").Value = "  +10.80%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "'149.19"
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").Value = "'18.25"
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("E34").Value = "  +5.33%  "
$ws.Range("D35").Value = "'5.20"
$ws.Range("E35").Value = "  +4.26%  "
$ws.Range("E36").Value = "  +8.46%  "
$ws.Range("E37").Value = "  +5.98%  "
$ws.Range("D38").Value = "'0.860"
$ws.Range("E38").Value = "  +7.70%  "
$ws.Range("D39").Value = "'34.12"
$ws.Range("E39").Value = "  +3.98%  "
$ws.Range("D40").Value = "'3.51"
$ws.Range("E40").Value = "  +8.82%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").Value = "'0.0555"
$ws.Range("E42").Value = "  +6.48%  "
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("E44").Value = "  +7.82%  "
$ws.Range("E45").Value = "  +15.31%  "
$ws.Range("D46").Value = "'0.0925"
$ws.Range("E46").Value = "  +5.33%  "
$ws.Range("D47").Value = "'257.87"
$ws.Range("E47").Value = "  +13.38%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0228"
$ws.Range("E48").Value = "  +5.36%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'10.18"
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("D50").Value = "'17.53"
$ws.Range("E50").Value = "  +6.35%  "
$ws.Range("D51").Value = "1.878.23"
$ws.Range("E51").Value = "  -2.94%  "
